$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (First VFrame, Last VFrame) before the existing
# "File Path" column (F), pushing it to H.
$ws.Columns("F:G").Insert()

# New header cells
$ws.Range("F1").Value = "First VFrame"
$ws.Range("G1").Value = "Last VFrame"

# First/Last VFrame numeric values per trial row
$firstFrames = @(415, 503, 376, 232, 345, 98)
$lastFrames  = @(477, 578, 452, 295, 410, 165)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $firstFrames[$i]
    $ws.Cells.Item($row, 7).Value = $lastFrames[$i]
}

# Fix up the selection to span the new used range
$ws.Range("A1:H7").Select()
